$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.505.29"
$ws.Range("E2").Value = "  +0.99%  "

$ws.Range("D3").Value = "2.240.39"
$ws.Range("E3").Value = "  -0.21%  "

$ws.Range("E4").Value = "  +0.33%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.82"
$ws.Range("E5").Value = "  -1.26%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.628"
$ws.Range("E6").Value = "  -0.79%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "75.13"
$ws.Range("E7").Value = "  -1.98%  "

$ws.Range("E8").Value = "  +0.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.621"
$ws.Range("E9").Value = "  -1.73%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.51"
$ws.Range("E10").Value = "  +5.38%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0947"
$ws.Range("E11").Value = "  -1.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.18"
$ws.Range("E12").Value = "  -0.40%  "

$ws.Range("E13").Value = "  +0.55%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.53"
$ws.Range("E14").Value = "  -2.00%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.858"
$ws.Range("E15").Value = "  -0.58%  "

$ws.Range("D16").Value = "2.247.94"
$ws.Range("E16").Value = "  -0.39%  "

$ws.Range("D17").Value = "42.290.00"
$ws.Range("E17").Value = "  +0.65%  "

$ws.Range("E18").Value = "  +4.56%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.18"
$ws.Range("E19").Value = "  +0.73%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.94"
$ws.Range("E20").Value = "  +0.05%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.57"
$ws.Range("E21").Value = "  +45.93%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "230.30"
$ws.Range("E22").Value = "  -0.61%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.19"
$ws.Range("E23").Value = "  -5.38%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.65"
$ws.Range("E24").Value = "  +2.22%  "

$ws.Range("E25").Value = "  +0.01%  "

$ws.Range("E26").Value = "  -1.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.30"
$ws.Range("E27").Value = "  -0.52%  "

$ws.Range("E28").Value = "  +4.51%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.91"
$ws.Range("E29").Value = "  -1.41%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.68"
$ws.Range("E30").Value = "  +0.55%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.95"
$ws.Range("E31").Value = "  +21.08%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0815"
$ws.Range("E32").Value = "  -2.01%  "

$ws.Range("E33").Value = "  -2.51%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.10"
$ws.Range("E34").Value = "  -11.07%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.125"
$ws.Range("E35").Value = "  -0.30%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.63"
$ws.Range("E36").Value = "  +1.40%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0313"
$ws.Range("E37").Value = "  +3.78%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "13.42"
$ws.Range("E38").Value = "  -5.87%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.17"
$ws.Range("E39").Value = "  -1.07%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.69"
$ws.Range("E40").Value = "  -4.35%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "63.45"
$ws.Range("E41").Value = "  +3.76%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.201"
$ws.Range("E42").Value = "  -0.89%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "106.46"
$ws.Range("E43").Value = "  -6.26%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.83"
$ws.Range("E44").Value = "  +1.39%  "

$ws.Range("E45").Value = "  +1.88%  "

$ws.Range("E46").Value = "  -0.06%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.41"
$ws.Range("E47").Value = "  +5.97%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.14"
$ws.Range("E48").Value = "  +0.04%  "

$ws.Range("E49").Value = "  +0.53%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.15"
$ws.Range("E50").Value = "  -0.86%  "

$ws.Range("E51").Value = "  +1.41%  "
